# Apply corrected data values to Sheet1 as described in the commit:
#   "corrected data added a brief sulfur analysis script with
#    stoichiometric ratios for NO3- and SO4(2-) reactions and
#    improved the plots (to be finished)."
#
# Column layout: A=sample, ..., M=weight, N=facies, O=TOC

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 2: fill in previously empty TOC value
$ws.Range("O2").Value = 1.2

# Row 7: correct TOC value
$ws.Range("O7").Value = 7.5

# Row 8: correct weight value
$ws.Range("M8").Value = 3.42

# Row 9: correct weight value, fill in TOC value
$ws.Range("M9").Value = 3.51
$ws.Range("O9").Value = 42.6

# Row 10: correct TOC value
$ws.Range("O10").Value = 3.67

# Row 15: correct weight value, fill in TOC value
$ws.Range("M15").Value = 3.34
$ws.Range("O15").Value = 44.0

# Row 17: fill in previously empty TOC value
$ws.Range("O17").Value = 20.05

$wb.Save()
